$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 500
$ws.Range("J32").Value = 500
$ws.Range("L32").Value = 500
$ws.Range("N32").Value = -1152
$ws.Range("H98").Value = 1011.15
$ws.Range("I98").Value = 652
$ws.Range("J98").Value = 1849.1666
$ws.Range("K98").Value = 652
$ws.Range("L98").Value = 1849.1666
$ws.Range("M98").Value = 846
$ws.Range("N98").Value = -4845.1666
$ws.Range("H122").Value = 1011.15
$ws.Range("I122").Value = 652
$ws.Range("J122").Value = 1849.1666
$ws.Range("K122").Value = 1956
$ws.Range("L122").Value = 5547.4998
$ws.Range("M122").Value = 494
$ws.Range("N122").Value = -10447.4998
$ws.Range("H129").Value = 1559.25
$ws.Range("J129").Value = 1919.25
$ws.Range("L129").Value = 5757.75
$ws.Range("N129").Value = -15757.75
$ws.Range("H135").Value = 6697.6924
$ws.Range("I135").Value = 5940.0835
$ws.Range("J135").Value = 15789
$ws.Range("K135").Value = 53460.7515
$ws.Range("L135").Value = 142101
$ws.Range("M135").Value = -50925.7515
$ws.Range("N135").Value = -147171

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 59074.168
$ws.Range("J23").Value = 38142.332
$ws.Range("L23").Value = 38142.332
$ws.Range("N23").Value = -38660.332
$ws.Range("H44").Value = 23625
$ws.Range("I44").Value = 20000
$ws.Range("J44").Value = 24833.334
$ws.Range("K44").Value = 20000
$ws.Range("L44").Value = 24833.334
$ws.Range("M44").Value = -19512
$ws.Range("N44").Value = -25809.334
$ws.Range("H45").Value = 191823.36
$ws.Range("I45").Value = 263194.62
$ws.Range("K45").Value = 263194.62
$ws.Range("M45").Value = -262817.62
$ws.Range("H55").Value = 21744.75
$ws.Range("J55").Value = 21744.75
$ws.Range("L55").Value = 21744.75
$ws.Range("N55").Value = -22374.75
$ws.Range("H61").Value = 242940.45
$ws.Range("I61").Value = 5716.731
$ws.Range("J61").Value = 628429
$ws.Range("K61").Value = 5716.731
$ws.Range("L61").Value = 628429
$ws.Range("M61").Value = -5504.731
$ws.Range("N61").Value = -628853
$ws.Range("H74").Value = 11112617
$ws.Range("I74").Value = 1264.3704
$ws.Range("J74").Value = 27779646
$ws.Range("K74").Value = 1264.3704
$ws.Range("L74").Value = 27779646
$ws.Range("M74").Value = -390.3704
$ws.Range("N74").Value = -27781394
$ws.Range("H77").Value = 11112617
$ws.Range("I77").Value = 1264.3704
$ws.Range("J77").Value = 27779646
$ws.Range("K77").Value = 6321.852
$ws.Range("L77").Value = 138898230
$ws.Range("M77").Value = -1953.852
$ws.Range("N77").Value = -138906966
$ws.Range("H80").Value = 42110
$ws.Range("J80").Value = 42110
$ws.Range("L80").Value = 42110
$ws.Range("N80").Value = -44106
$ws.Range("H83").Value = 42110
$ws.Range("J83").Value = 42110
$ws.Range("L83").Value = 126330
$ws.Range("N83").Value = -136314
$ws.Range("H136").Value = 242940.45
$ws.Range("I136").Value = 5716.731
$ws.Range("J136").Value = 628429
$ws.Range("K136").Value = 17150.193
$ws.Range("L136").Value = 1885287
$ws.Range("M136").Value = -14600.193
$ws.Range("N136").Value = -1890387

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 216647
$ws.Range("I107").Value = 275268.9
$ws.Range("K107").Value = 275268.9
$ws.Range("M107").Value = -273348.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3842117.5
$ws.Range("I58").Value = 5953698.5
$ws.Range("J58").Value = 557435.5600000001
$ws.Range("K58").Value = 5953698.5
$ws.Range("L58").Value = 557435.5600000001
$ws.Range("M58").Value = -5953495.5
$ws.Range("N58").Value = -557841.5600000001
$ws.Range("H136").Value = 3842117.5
$ws.Range("I136").Value = 5953698.5
$ws.Range("J136").Value = 557435.5600000001
$ws.Range("K136").Value = 17861095.5
$ws.Range("L136").Value = 1672306.68
$ws.Range("M136").Value = -17858545.5
$ws.Range("N136").Value = -1677406.68
$ws.Range("H140").Value = 24418.182
$ws.Range("I140").Value = 5000
$ws.Range("J140").Value = 25025
$ws.Range("K140").Value = 5000
$ws.Range("L140").Value = 25025
$ws.Range("M140").Value = 180
$ws.Range("N140").Value = -35385

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 737.0645
$ws.Range("I121").Value = 367.125
$ws.Range("J121").Value = 865.73914
$ws.Range("K121").Value = 1101.375
$ws.Range("L121").Value = 2597.21742
$ws.Range("M121").Value = 208.625
$ws.Range("N121").Value = -5217.21742

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2766.724
$ws.Range("I102").Value = 2628.8635
$ws.Range("J102").Value = 3200
$ws.Range("K102").Value = 2628.8635
$ws.Range("L102").Value = 3200
$ws.Range("M102").Value = -1006.8635
$ws.Range("N102").Value = -6444
$ws.Range("H122").Value = 64993576
$ws.Range("I122").Value = 88736344
$ws.Range("J122").Value = 33336556
$ws.Range("K122").Value = 266209032
$ws.Range("L122").Value = 100009668
$ws.Range("M122").Value = -266206582
$ws.Range("N122").Value = -100014568
$ws.Range("H126").Value = 17658.334
$ws.Range("I126").Value = 25112.5
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 75337.5
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = -72867.5
$ws.Range("N126").Value = -13190
$ws.Range("H132").Value = 3147783.5
$ws.Range("I132").Value = 4764652.5
$ws.Range("K132").Value = 14293957.5
$ws.Range("M132").Value = -14291427.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12505801
$ws.Range("I132").Value = 13895002
$ws.Range("J132").Value = 2992.5
$ws.Range("K132").Value = 41685006
$ws.Range("L132").Value = 8977.5
$ws.Range("M132").Value = -41682476
$ws.Range("N132").Value = -14037.5
$ws.Range("H136").Value = 20542.715
$ws.Range("I136").Value = 17574.875
$ws.Range("J136").Value = 24499.834
$ws.Range("K136").Value = 52724.625
$ws.Range("L136").Value = 73499.50199999999
$ws.Range("M136").Value = -50174.625
$ws.Range("N136").Value = -78599.50199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 42110
$ws.Range("J64").Value = 42110
$ws.Range("L64").Value = 42110
$ws.Range("N64").Value = -42606
$ws.Range("H67").Value = 42110
$ws.Range("J67").Value = 42110
$ws.Range("L67").Value = 42110
$ws.Range("N67").Value = -43826
